# imprimindo no main linha a linha o json[balanco_detalhado][subsistemas][fontes]
#
# A linha 87 (ultima linha de dados existente) e replicada ate a linha 108,
# simulando a impressao, linha a linha, das fontes de cada subsistema.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BalancoResumido")

# Linha-modelo com os dados (colunas A:Q) a serem repetidos
$template = $ws.Range("A87:Q87")

for ($row = 88; $row -le 108; $row++) {
    $dest = $ws.Range("A" + $row + ":Q" + $row)
    $template.Copy($dest)
}

# A antiga linha 87 tinha uma celula vazia extra em R87; ela deixa de existir
# e a ultima linha da nova sequencia (108) passa a receber essa marca vazia.
$ws.Range("R87").Copy($ws.Range("R108"))
$ws.Range("R87").ClearContents()
